$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header cells in row 1 (T1:AC1): cost-related columns ---
$ws.Range("T1").Value  = "LA.Avg.Gross.Weekly.Cost.Per.Person"
$ws.Range("U1").Value  = "LA.Cost.Of.Care.18-64"
$ws.Range("V1").Value  = "LA.Cost.Of.Care.65-74"
$ws.Range("W1").Value  = "LA.Cost.Of.Care.75-84"
$ws.Range("X1").Value  = "LA.Cost.Of.Care.75pl"
$ws.Range("Y1").Value  = "Other Supplementary LA Cost Metrics"
$ws.Range("Z1").Value  = "…"
$ws.Range("AA1").Value = ".."
$ws.Range("AB1").Value = ".."
$ws.Range("AC1").Value = ".."

# --- Row 3: "Res / Nursing" annotations over the occupancy + new cost cols ---
$ws.Range("N3").Value = "Res / Nursing"
$ws.Range("O3").Value = "Res / Nursing"
$ws.Range("P3").Value = "Res / Nursing"
$ws.Range("Q3").Value = "Res / Nursing"
$ws.Range("R3").Value = "Res / Nursing?"
$ws.Range("S3").Value = "Res / Nursing?"
$ws.Range("T3").Value = "Res / Nursing"
$ws.Range("U3").Value = "Res / Nursing"
$ws.Range("V3").Value = "Res / Nursing"
$ws.Range("W3").Value = "Res / Nursing"
$ws.Range("X3").Value = "Res / Nursing"
$ws.Range("Y3").Value = "TBD!"

# --- Row 5: explanatory / scoping notes ---
$ws.Range("A5").Value = "~150 LAs"
$ws.Range("C5").Value = "A row for every year between now and 2037"
$ws.Range("F5").Value = "From ONS SRC dataset"
$ws.Range("G5").Value = "From ONS SRC dataset"
$ws.Range("H5").Value = "From ONS SRC dataset"
$ws.Range("I5").Value = "From ONS SRC dataset"
$ws.Range("J5").Value = "From ONS SRC dataset"
$ws.Range("K5").Value = "From ONS SRC dataset"
$ws.Range("L5").Value = "From ONS SRC dataset"
$ws.Range("M5").Value = "From ONS SRC dataset"
$ws.Range("U5").Value = "Estiamted LA-funded occupancy, x LA Avg Gross Weekly Cost of Care"
$ws.Range("V5").Value = "Estiamted LA-funded occupancy, x LA Avg Gross Weekly Cost of Care"
$ws.Range("W5").Value = "Estiamted LA-funded occupancy, x LA Avg Gross Weekly Cost of Care"
$ws.Range("X5").Value = "Estiamted LA-funded occupancy, x LA Avg Gross Weekly Cost of Care"

# --- Column widths for the newly added columns ---
$ws.Range("U1").ColumnWidth = 21.28515625
$ws.Range("V1").ColumnWidth = 21.28515625
$ws.Range("W1").ColumnWidth = 21.28515625
$ws.Range("X1").ColumnWidth = 21.28515625
$ws.Range("Y1").ColumnWidth = 33.28515625

# --- Scroll the view so column L is leftmost, and select N5 (matches author's end state) ---
$excel.ActiveWindow.ScrollColumn = 12
$ws.Range("N5").Select()
